$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '54.758.08'

$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +6.21%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.427.73'

$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +7.43%  '

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.17%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '478.00'

$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +11.30%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '138.95'

$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +21.20%  '

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.08%  '

$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +11.07%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.454.80'

$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +8.41%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0957'

$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +15.64%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.47'

$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +7.05%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.321'

$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +10.42%  '

$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +2.23%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.849.16'

$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +7.95%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '54.917.20'

$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +6.46%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.45'

$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +13.25%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000134'

$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +20.05%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.449.46'

$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +8.31%  '

$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +13.39%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '9.82'

$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +17.57%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '312.03'

$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +8.17%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.994'

$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.17%  '

$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +16.36%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '57.10'

$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +9.06%  '

$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.13%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.401'

$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +13.05%  '

$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +20.82%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.551.15'

$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +9.34%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.33'

$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +11.52%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0766'

$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +25.23%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.997'

$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.06%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '148.37'

$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +3.81%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '17.87'

$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +10.70%  '

$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +14.82%  '

$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +14.21%  '

$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +18.21%  '

$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +11.06%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.59'

$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +11.27%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '33.38'

$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +6.07%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.993'

$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.38%  '

$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +10.15%  '

$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +13.79%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0541'

$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +12.54%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.28'

$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +18.04%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.13'

$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.18%  '

$ws.Range('B46').Value = 'RenderToken'

$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.64'

$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +21.62%  '

$ws.Range('B47').Value = 'Bittensor'

$ws.Range('C47').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '254.41'

$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +35.31%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0886'

$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +12.99%  '

$ws.Range('B49').Value = 'VeChain'

$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0221'

$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +13.23%  '

$ws.Range('B50').Value = 'Maker'

$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.916.78'

$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +4.53%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '17.01'

$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +12.99%  '
